$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.245.79'
$ws.Range("E2").Value = '  +1.95%  '

$ws.Range("D3").Value = '2.304.63'
$ws.Range("E3").Value = '  +1.29%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.25%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.503'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.89%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.517'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +5.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.46'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +10.38%  '

$ws.Range("E11").Value = '  +0.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.77'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +13.11%  '

$ws.Range("E13").Value = '  +1.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.04%  '

$ws.Range("D15").Value = '2.664.43'
$ws.Range("E15").Value = '  +1.15%  '

$ws.Range("D16").Value = '2.348.39'
$ws.Range("E16").Value = '  +4.26%  '

$ws.Range("E17").Value = '  +1.91%  '

$ws.Range("D18").Value = '43.083.58'
$ws.Range("E18").Value = '  +1.70%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.68'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +11.63%  '

$ws.Range("E20").Value = '  +4.54%  '

$ws.Range("D21").Value = '0.0₃0909'
$ws.Range("E21").Value = '  +2.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.64'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.53%  '

$ws.Range("E24").Value = '  +15.61%  '

$ws.Range("E25").Value = '  +0.43%  '

$ws.Range("E26").Value = '  +0.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '34.97'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '168.06'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.05%  '

$ws.Range("E31").Value = '  +1.78%  '

$ws.Range("E32").Value = '  -0.09%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.05'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.22%  '

$ws.Range("E35").Value = '  -0.22%  '

$ws.Range("E36").Value = '  +0.24%  '

$ws.Range("E37").Value = '  +1.63%  '

$ws.Range("E38").Value = '  +2.30%  '

$ws.Range("E39").Value = '  +4.57%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.102'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.68%  '

$ws.Range("E41").Value = '  +0.89%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.34'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.80%  '

$ws.Range("D43").Value = '1.992.99'
$ws.Range("E43").Value = '  +2.20%  '

$ws.Range("E44").Value = '  +5.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.65%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '17.78'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.01%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.88'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.10%  '

$ws.Range("E49").Value = '  +6.62%  '

$ws.Range("D50").Value = '2.530.47'
$ws.Range("E50").Value = '  +0.94%  '

# Row 51: coin changed from THORChain to BitcoinSV
$ws.Range("B51").Value = 'BitcoinSV'
$ws.Range("C51").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '70.94'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.28%  '
